# Refresh the cryptos list (Price / Volume(1h) columns), as produced by the
# "Updated cryptos list ... with GitHub Actions" scheduled job. Two rows
# (Toncoin/RenderToken and WEMIXToken/ApeXProtocol) also swapped rank
# position, so their Coin name + Link cells are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their text formatting so that
# numeric-looking strings (e.g. "71.20") are not coerced into numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "52.489.30"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "2.928.34"
$ws.Range("E3").Value = "  +4.49%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "353.14"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "113.04"
$ws.Range("E6").Value = "  +3.61%  "
$ws.Range("D7").Value = "0.563"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "40.25"
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").Value = "0.0865"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "20.19"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "7.86"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "3.386.76"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "2.927.22"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "0.993"
$ws.Range("E17").Value = "  +5.06%  "
$ws.Range("D18").Value = "52.494.34"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +5.27%  "
$ws.Range("D21").Value = "14.44"
$ws.Range("E21").Value = "  +6.72%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").Value = "71.20"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "271.67"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").Value = "27.08"
$ws.Range("E26").Value = "  +3.89%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "10.66"
$ws.Range("E29").Value = "  +3.10%  "
$ws.Range("D30").Value = "38.23"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").Value = "6.55"
$ws.Range("E31").Value = "  +5.81%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "6.21"
$ws.Range("E32").Value = "  +8.92%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.25"
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "53.30"
$ws.Range("E34").Value = "  +2.60%  "
$ws.Range("D35").Value = "0.0941"
$ws.Range("E35").Value = "  +10.04%  "
$ws.Range("D36").Value = "0.0454"
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "3.36"
$ws.Range("E38").Value = "  +7.31%  "
$ws.Range("E39").Value = "  +4.78%  "
$ws.Range("D40").Value = "18.81"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "2.74"
$ws.Range("E41").Value = "  +10.37%  "
$ws.Range("D42").Value = "24.64"
$ws.Range("E42").Value = "  +13.21%  "
$ws.Range("E43").Value = "  +2.00%  "
$ws.Range("D44").Value = "122.95"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.61"
$ws.Range("E45").Value = "  +8.15%  "
$ws.Range("B46").Value = "WEMIXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "2.20"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "2.220.26"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("D48").Value = "3.57"
$ws.Range("D49").Value = "0.264"
$ws.Range("E49").Value = "  +24.57%  "
$ws.Range("D50").Value = "0.0341"
$ws.Range("E50").Value = "  +16.04%  "
$ws.Range("D51").Value = "0.959"
$ws.Range("E51").Value = "  +5.34%  "

# Restore the default cell style (no explicit NumberFormat) to match
# the original workbook formatting.
$dataRange.Style = "Normal"
